$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 75
$ws.Range("A75").Value2 = 111949317
$ws.Range("B75").Value2 = 96265
$ws.Range("D75").Value2 = "LC"
$ws.Range("E75").Value2 = 219790
$ws.Range("F75").Value2 = "Fläcknycklar"
$ws.Range("G75").Value2 = "Dactylorhiza maculata"
$ws.Range("H75").Value2 = "(L.) Soó"
$ws.Range("I75").NumberFormat = "@"
$ws.Range("I75").Value2 = ""
$ws.Range("I75").Style = "Normal"
$ws.Range("J75").ClearContents()
$ws.Range("Q75").Value2 = 580500.003505226
$ws.Range("R75").Value2 = 7053328.641698814
$ws.Range("S75").Value2 = 2
$ws.Range("Z75").Value2 = "18:54"
$ws.Range("AB75").Value2 = "18:54"
$ws.Range("AW75").Value2 = "Kim Hultgren"
$ws.Range("AX75").Value2 = "Kim Hultgren"

# Row 76
$ws.Range("A76").Value2 = 111950184
$ws.Range("B76").Value2 = 56543
$ws.Range("D76").Value2 = "NT"
$ws.Range("E76").Value2 = 103021
$ws.Range("F76").Value2 = "Talltita"
$ws.Range("G76").Value2 = "Poecile montanus"
$ws.Range("H76").Value2 = "(Conrad von Baldenstein, 1827)"
$ws.Range("I76").NumberFormat = "@"
$ws.Range("I76").Value2 = ""
$ws.Range("I76").Style = "Normal"
$ws.Range("Q76").Value2 = 580446.7330953531
$ws.Range("R76").Value2 = 7053301.910512885
$ws.Range("S76").Value2 = 10
$ws.Range("Z76").Value2 = "19:37"
$ws.Range("AB76").Value2 = "19:37"

# Row 77
$ws.Range("A77").Value2 = 111949678
$ws.Range("B77").Value2 = 96348
$ws.Range("D77").Value2 = "VU"
$ws.Range("E77").Value2 = 220787
$ws.Range("F77").Value2 = "Knärot"
$ws.Range("G77").Value2 = "Goodyera repens"
$ws.Range("H77").Value2 = "(L.) R. Br."
$ws.Range("I77").NumberFormat = "@"
$ws.Range("I77").Value2 = "7"
$ws.Range("I77").Style = "Normal"
$ws.Range("Q77").Value2 = 580467.4207067642
$ws.Range("R77").Value2 = 7053330.04139028
$ws.Range("S77").Value2 = 2
$ws.Range("Z77").Value2 = "19:11"
$ws.Range("AB77").Value2 = "19:11"

# Row 78
$ws.Range("A78").Value2 = 111949575
$ws.Range("B78").Value2 = 96348
$ws.Range("D78").Value2 = "VU"
$ws.Range("E78").Value2 = 220787
$ws.Range("F78").Value2 = "Knärot"
$ws.Range("G78").Value2 = "Goodyera repens"
$ws.Range("H78").Value2 = "(L.) R. Br."
$ws.Range("I78").NumberFormat = "@"
$ws.Range("I78").Value2 = "15"
$ws.Range("I78").Style = "Normal"
$ws.Range("J78").Value2 = "plantor/tuvor"
$ws.Range("Q78").Value2 = 580471.3517951096
$ws.Range("R78").Value2 = 7053333.257918903
$ws.Range("S78").Value2 = 1
$ws.Range("Z78").Value2 = "19:05"
$ws.Range("AB78").Value2 = "19:05"
$ws.Range("AW78").Value2 = "Kamilla Andersson"
$ws.Range("AX78").Value2 = "Kamilla Andersson"

# Row 79
$ws.Range("A79").Value2 = 111949591
$ws.Range("I79").NumberFormat = "@"
$ws.Range("I79").Value2 = "20"
$ws.Range("I79").Style = "Normal"
$ws.Range("Q79").Value2 = 580476.1122211452
$ws.Range("R79").Value2 = 7053321.356648902
$ws.Range("Z79").Value2 = "19:07"
$ws.Range("AB79").Value2 = "19:07"
